# Refreshes the cryptocurrency snapshot (price + 1h volume%) for the rows
# that moved, and fixes three mis-ranked coin entries (Chainlink/ShibaInu
# swap at rows 19-20, Quant/mCoin swap at rows 48-49, Cronos -> BitcoinSV
# at row 51) to match the refreshed coinranking.com leaderboard.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, already-dotted/percent strings -
# Excel will not misparse these as numbers).
$textUpdates = [ordered]@{
    'D2' = '28.529.21'
    'E2' = '  +3.82%  '
    'D3' = '1.590.76'
    'E3' = '  +1.28%  '
    'E4' = '  +0.85%  '
    'E5' = '  +0.75%  '
    'E6' = '  -0.03%  '
    'E7' = '  +0.79%  '
    'E8' = '  +7.37%  '
    'E9' = '  +0.25%  '
    'E10' = '  +0.36%  '
    'E11' = '  +1.74%  '
    'D12' = '1.819.08'
    'E12' = '  +1.37%  '
    'D13' = '1.590.18'
    'E13' = '  +1.36%  '
    'E14' = '  +2.00%  '
    'E15' = '  -0.37%  '
    'D16' = '28.564.19'
    'E16' = '  +3.96%  '
    'E17' = '  +1.14%  '
    'E18' = '  +2.85%  '
    'B19' = 'Chainlink'
    'C19' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'E19' = '  -0.58%  '
    'B20' = 'ShibaInu'
    'C20' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D20' = '0.0₃0707'
    'E20' = '  -0.02%  '
    'E21' = '  +0.79%  '
    'E22' = '  -1.89%  '
    'E23' = '  -0.80%  '
    'E24' = '  +1.38%  '
    'E25' = '  +0.77%  '
    'E26' = '  +0.85%  '
    'E27' = '  -0.60%  '
    'E28' = '  -0.41%  '
    'E29' = '  +0.82%  '
    'E30' = '  -0.52%  '
    'E31' = '  +0.05%  '
    'E32' = '  +0.51%  '
    'D34' = '1.403.81'
    'E34' = '  -3.59%  '
    'E35' = '  -0.89%  '
    'E36' = '  -10.22%  '
    'E37' = '  +0.95%  '
    'E38' = '  +11.47%  '
    'E39' = '  -0.62%  '
    'E40' = '  +0.75%  '
    'E41' = '  +0.25%  '
    'E42' = '  +0.78%  '
    'E43' = '  -0.58%  '
    'E44' = '  +0.48%  '
    'E45' = '  +0.74%  '
    'E46' = '  -1.65%  '
    'D47' = '1.728.41'
    'E47' = '  +1.28%  '
    'B48' = 'Quant'
    'C48' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E48' = '  +0.53%  '
    'B49' = 'mCoin'
    'C49' = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
    'E49' = '  +1.26%  '
    'E50' = '  +0.97%  '
    'B51' = 'BitcoinSV'
    'C51' = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
    'E51' = '  +12.76%  '
}

# Number-looking price strings (e.g. "212.90", "41.50"). These need the
# Text-format round-trip below so Excel keeps the literal string (with its
# trailing zeros / dot-grouping) instead of silently coercing to a double.
$numericLookingUpdates = [ordered]@{
    'D5' = '212.90'
    'D8' = '24.41'
    'D9' = '0.251'
    'D11' = '0.0887'
    'D14' = '0.529'
    'D17' = '63.17'
    'D18' = '232.72'
    'D19' = '7.49'
    'D22' = '4.04'
    'D23' = '9.34'
    'D24' = '1.98'
    'D26' = '15.27'
    'D38' = '2.63'
    'D40' = '0.543'
    'D45' = '0.982'
    'D46' = '63.14'
    'D48' = '87.33'
    'D49' = '2.13'
    'D51' = '41.50'
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
    # Restore the default "Normal" style so no stray formatting is left
    # behind on the cell once the literal text is safely stored.
    $cell.Style = "Normal"
}
